$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.259.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.242.75"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.55%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.77"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.13"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.04%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0988"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.07"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "36.45"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +10.57%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.76"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.576.01"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.97"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -7.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.865"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.246.54"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.135.02"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0969"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.20"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +9.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.66"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.46"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.54%  "

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.99"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.06"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.54"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.58%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.16%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0720"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.07%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.30"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.69"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.84"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.79"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +20.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0281"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.30"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.43%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "67.51"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.90%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.89"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.33"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.95"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -12.21%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.53%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.190"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +10.37%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.19"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.34"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.89%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Celestia"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.05"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.69%  "
